$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the existing 5 match rows (currently rows 2-6) again below the
# current data, growing the sheet from A1:K6 to A1:K11. The appended block
# repeats the same 5 rows, but with the "Sharjah"/"Abu Dhabi" rows (4 and 5)
# swapped in order: row7<-row2, row8<-row3, row9<-row5, row10<-row4, row11<-row6.
# Use Copy/PasteSpecial (rather than .Value=) so the destination cells keep
# the exact same "number stored as text" cell type as the source cells.

$sourceRows = @(2, 3, 5, 4, 6)
$destRow = 7

foreach ($srcRow in $sourceRows) {
    $ws.Range("A$srcRow`:K$srcRow").Copy()
    $ws.Range("A$destRow").PasteSpecial()
    $destRow++
}
